# edit.ps1 — applies the "Update Danh sach de thi cong khai" change:
#   * Deletes the whole "Thiet ke lai Trang Bang dieu khien (Dashboard
#     Redesign)" bullet block (heading + Hien tai / Y tuong / 3 sub-bullets /
#     Loi ich) from section "A. Cai thien Bo cuc va Luong lam viec", so the
#     two remaining blocks ("Cai thien 'Trang thai rong'" and "Bo cuc 2 cot
#     cho Trang Tim kiem") shift up.
#   * Moves the hidden "_GoBack" bookmark from its old spot (end of section
#     "C", right before the page break) to the start of the new first
#     heading of that block ("Cai thien 'Trang thai rong' ...").
#   * Drops the stale lastRenderedPageBreak cache marker that is left
#     behind on the "C. Nang cao Trai nghiem Cot loi" heading now that the
#     page no longer breaks there.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Delete the entire "Dashboard Redesign" bullet block.
# ---------------------------------------------------------------------
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($startPara -eq $null -and $p.Range.Text -like "Thiết kế lại Trang Bảng điều khiển*") {
        $startPara = $p
    }
    if ($startPara -ne $null -and $endPara -eq $null -and $p.Range.Text -like "*Làm cho Dashboard trở nên hữu ích*") {
        $endPara = $p
    }
}
if ($startPara -ne $null -and $endPara -ne $null) {
    $blockRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $blockRange.Delete()
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark onto the new first heading of the block
#    ("Cải thiện "Trạng thái rỗng" (Better Empty States):"). Adding a
#    bookmark with a name that already exists relocates it, so the stale
#    one near "C. Nâng cao ..." disappears automatically.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Cải thiện*Trạng thái rỗng*Better Empty States*") {
        $pos = $p.Range.Start
        $bmRange = $d.Range($pos, $pos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}

# ---------------------------------------------------------------------
# 3) Clear the leftover lastRenderedPageBreak marker on "C. Nâng cao
#    Trải nghiệm Cốt lõi" — with less content above it, the page no
#    longer breaks right there. Touching the run (append + remove one
#    character right at its end) forces it to be rebuilt without the
#    stale cached marker while keeping its text/formatting identical.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "C. Nâng cao Trải nghiệm Cốt lõi*") {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $endPos = $r.End
        $touch = $d.Range($endPos, $endPos)
        $touch.InsertAfter("Z")
        $cleanup = $d.Range($endPos, $endPos + 1)
        $cleanup.Delete()
        break
    }
}
